$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------------
# Sheet "Overview": the two file rows (9675f9d4...md and aa0ab8e1...md) swap
# places - aa0ab8e1 moves up to row 2 (keeps "Handed back" status) and
# 9675f9d4 moves down to row 3 and gets a new "Ready for handoff" status.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8e57aed0dc9a09ccbd02ccc8da1353424e3d70ac/e2e/aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md", $missing, $missing, "aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8e57aed0dc9a09ccbd02ccc8da1353424e3d70ac/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", $missing, $missing, "9675f9d4-d77a-4429-af2a-e4d43d867617.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8e57aed0dc9a09ccbd02ccc8da1353424e3d70ac/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row swap as above, with the 9675f9d4 row (now row 3)
# getting status "Ready for handoff" and a new, later handoff datetime.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "2016-01-25 07:54:50"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8e57aed0dc9a09ccbd02ccc8da1353424e3d70ac/e2e/aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md", $missing, $missing, "aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec71c7993df7839c1f2207009350ac9f937abf71/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/aa0ab8e1-ea6a-411d-9970-d3af90ea867b.0e65d7b5f42d83ab517e0727b1b5c5dcefb556f5.zh-cn.xlf", $missing, $missing, "aa0ab8e1-ea6a-411d-9970-d3af90ea867b.0e65d7b5f42d83ab517e0727b1b5c5dcefb556f5.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ae1270bcdc8534d5339decd743c5ffb2275a2dab/e2e/aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md", $missing, $missing, "aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5a2af4fe5a044a6c77ce8e320d4ae60d91aa35a8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/aa0ab8e1-ea6a-411d-9970-d3af90ea867b.0e65d7b5f42d83ab517e0727b1b5c5dcefb556f5.zh-cn.xlf", $missing, $missing, "aa0ab8e1-ea6a-411d-9970-d3af90ea867b.0e65d7b5f42d83ab517e0727b1b5c5dcefb556f5.zh-cn.xlf") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8e57aed0dc9a09ccbd02ccc8da1353424e3d70ac/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", $missing, $missing, "9675f9d4-d77a-4429-af2a-e4d43d867617.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec71c7993df7839c1f2207009350ac9f937abf71/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.zh-cn.xlf", $missing, $missing, "9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ae1270bcdc8534d5339decd743c5ffb2275a2dab/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", $missing, $missing, "9675f9d4-d77a-4429-af2a-e4d43d867617.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5a2af4fe5a044a6c77ce8e320d4ae60d91aa35a8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.zh-cn.xlf", $missing, $missing, "9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.zh-cn.xlf") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8e57aed0dc9a09ccbd02ccc8da1353424e3d70ac/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": same row swap, with the 9675f9d4 row (now row 3) getting
# status "Ready for handoff" and a new, later handoff datetime.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "2016-01-25 07:55:06"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8e57aed0dc9a09ccbd02ccc8da1353424e3d70ac/e2e/aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md", $missing, $missing, "aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbab10a1e765cf6149063ad6cd37cd6dbc666827/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/aa0ab8e1-ea6a-411d-9970-d3af90ea867b.0e65d7b5f42d83ab517e0727b1b5c5dcefb556f5.de-de.xlf", $missing, $missing, "aa0ab8e1-ea6a-411d-9970-d3af90ea867b.0e65d7b5f42d83ab517e0727b1b5c5dcefb556f5.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/54d41df6777fcba3e49e891f507c67e262f885b1/e2e/aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md", $missing, $missing, "aa0ab8e1-ea6a-411d-9970-d3af90ea867b.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2266e5ebcb482b38b8f6a3beb114362c30ea5c88/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/aa0ab8e1-ea6a-411d-9970-d3af90ea867b.0e65d7b5f42d83ab517e0727b1b5c5dcefb556f5.de-de.xlf", $missing, $missing, "aa0ab8e1-ea6a-411d-9970-d3af90ea867b.0e65d7b5f42d83ab517e0727b1b5c5dcefb556f5.de-de.xlf") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8e57aed0dc9a09ccbd02ccc8da1353424e3d70ac/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", $missing, $missing, "9675f9d4-d77a-4429-af2a-e4d43d867617.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbab10a1e765cf6149063ad6cd37cd6dbc666827/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.de-de.xlf", $missing, $missing, "9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/54d41df6777fcba3e49e891f507c67e262f885b1/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", $missing, $missing, "9675f9d4-d77a-4429-af2a-e4d43d867617.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2266e5ebcb482b38b8f6a3beb114362c30ea5c88/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.de-de.xlf", $missing, $missing, "9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.de-de.xlf") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8e57aed0dc9a09ccbd02ccc8da1353424e3d70ac/.localization-config", $missing, $missing, ".localization-config") | Out-Null
